$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change mass for "Phased Array (Pivoting)" row (row 4) from 2 to 1
$ws.Range("G4").Value = 1

# Add a new "RELAY" range-indicator column (J/K) to the phased array rows
# J column: label "RELAY"; K column: maximum relay transmission range
# (values need to stay as text -- temporarily force Text format so Excel
#  doesn't coerce the scientific-notation-looking strings into numbers,
#  then strip the formatting back off so no stray style sticks around)
$ws.Range("K3:K7").NumberFormat = "@"

$ws.Range("J7").Value = "RELAY"
$ws.Range("K7").Value = "1.0e+11"

$ws.Range("J6").Value = "RELAY"
$ws.Range("K6").Value = "2.0e+11"

$ws.Range("J5").Value = "RELAY"
$ws.Range("K5").Value = "2.0e+11"

$ws.Range("J4").Value = "RELAY"
$ws.Range("K4").Value = "2.0e+11"

$ws.Range("J3").Value = "RELAY"
$ws.Range("K3").Value = "1.0e+12"

$ws.Range("K3:K7").ClearFormats()

$ws.Range("M12").Select()
